# Add a second row of data below the existing "test" entry in A1:
#   - A2 gets the new string value "test02" (creates a new shared-string entry)
#   - the sheet's used range (dimension) grows to A1:A2
#   - the active selection moves on to A3, as if the user had just typed
#     the value into A2 and pressed Enter
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test02"
$ws.Range("A3").Select() | Out-Null
